# chore: adapt column header formatting to respective input file names
#
# - rename the "_old" / "_new" header-row suffixes to "_FV2410" / "_FV2504"
# - turn the used range into a native Excel Table (Table1)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# --- 1. rewrite the header row (row 1) -------------------------------------
$headerMap = @{
    1  = "Segmentname_FV2410"
    2  = "Segmentgruppe_FV2410"
    3  = "Segment_FV2410"
    4  = "Datenelement_FV2410"
    5  = "Segment ID_FV2410"
    6  = "Code_FV2410"
    7  = "Qualifier_FV2410"
    8  = "Beschreibung_FV2410"
    9  = "Bedingungsausdruck_FV2410"
    10 = "Bedingung_FV2410"
    11 = "diff"
    12 = "Segmentname_FV2504"
    13 = "Segmentgruppe_FV2504"
    14 = "Segment_FV2504"
    15 = "Datenelement_FV2504"
    16 = "Segment ID_FV2504"
    17 = "Code_FV2504"
    18 = "Qualifier_FV2504"
    19 = "Beschreibung_FV2504"
    20 = "Bedingungsausdruck_FV2504"
    21 = "Bedingung_FV2504"
}

foreach ($col in $headerMap.Keys) {
    $ws.Cells.Item(1, $col).Value() = $headerMap[$col]
}

# --- 2. turn A1:U83 into an Excel Table (Table1) ----------------------------
$tableRange = $ws.Range("A1:U83")
$lo = $ws.ListObjects.Add(1, $tableRange, [System.Type]::Missing, 1, [System.Type]::Missing)
$lo.Name = "Table1"

# --- 3. freeze the header row -----------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
